$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.938.70"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "2.255.81"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'317.53"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'101.42"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.554"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'37.12"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "'7.65"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "2.601.78"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "'0.859"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "2.260.65"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "43.846.69"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").Value = "'13.51"
$ws.Range("E19").Value = "  -6.27%  "
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  +2.45%  "
$ws.Range("D21").Value = "'6.54"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'65.84"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "'3.11"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'235.09"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'10.15"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").Value = "'37.06"
$ws.Range("E29").Value = "  +3.63%  "
$ws.Range("D30").Value = "'6.22"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "'159.30"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").Value = "'20.17"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "'0.0850"
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("D34").Value = "'2.70"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").Value = "  +11.57%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D37").Value = "'3.07"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").Value = "'16.17"
$ws.Range("E39").Value = "  +20.20%  "
$ws.Range("D40").Value = "'3.71"
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").Value = "'4.21"
$ws.Range("E41").Value = "  -4.77%  "
$ws.Range("D42").Value = "'0.0315"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "1.815.33"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("D45").Value = "'75.76"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'82.40"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.197"
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("D48").Value = "'5.22"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").Value = "'104.99"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.68"
$ws.Range("E50").Value = "  +6.84%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'58.31"
$ws.Range("E51").Value = "  +0.11%  "
